# Generate Report for Archive
# - Status changes from "Ready for handoff" to "In Translation" for the two
#   localized files, on the Overview sheet (zh-cn/de-de status columns) as
#   well as on the per-locale "zh-cn" and "de-de" report sheets.
# - The "Status" column narrows to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) and de-de (col F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Column got narrower now that "In Translation" is shorter than
# "Ready for handoff". (The engine snaps ColumnWidth to the nearest
# renderable increment; 12.835 lands on the same rendered width as the
# narrower column from the authoritative edit.)
$newColumnWidth = 12.835
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# --- Per-locale sheets: "Status" is column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
